{"js": "// The document's sole table is a 20x5 grid of arithmetic equations\n// (\"a+b=\" / \"a-b=\"). The commit replaces the text of every cell with a\n// newly generated equation while leaving all formatting untouched.\n// Address each cell positionally (0-based row/col) and overwrite its\n// text via the Office.js TableCell.value setter so run/paragraph\n// formatting (font, size, justification) is preserved exactly.\nconst table = context.document.body.tables.getFirst();\n\ntable.getCell(0, 0).value = \"26+36=\";\ntable.getCell(0, 1).value = \"64+4=\";\ntable.getCell(0, 2).value = \"44+54=\";\ntable.getCell(0, 3).value = \"98-12=\";\ntable.getCell(0, 4).value = \"91-27=\";\ntable.getCell(1, 0).value = \"44+24=\";\ntable.getCell(1, 1).value = \"44-32=\";\ntable.getCell(1, 2).value = \"37+18=\";\ntable.getCell(1, 3).value = \"66+25=\";\ntable.getCell(1, 4).value = \"4+28=\";\ntable.getCell(2, 0).value = \"40+51=\";\ntable.getCell(2, 1).value = \"13+28=\";\ntable.getCell(2, 2).value = \"66+3=\";\ntable.getCell(2, 3).value = \"21+29=\";\ntable.getCell(2, 4).value = \"66-65=\";\ntable.getCell(3, 0).value = \"42-42=\";\ntable.getCell(3, 1).value = \"92-6=\";\ntable.getCell(3, 2).value = \"1+77=\";\ntable.getCell(3, 3).value = \"38+2=\";\ntable.getCell(3, 4).value = \"8+75=\";\ntable.getCell(4, 0).value = \"28-10=\";\ntable.getCell(4, 1).value = \"98-9=\";\ntable.getCell(4, 2).value = \"36+7=\";\ntable.getCell(4, 3).value = \"33-22=\";\ntable.getCell(4, 4).value = \"53+44=\";\ntable.getCell(5, 0).value = \"2+61=\";\ntable.getCell(5, 1).value = \"97-76=\";\ntable.getCell(5, 2).value = \"5+93=\";\ntable.getCell(5, 3).value = \"71-20=\";\ntable.getCell(5, 4).value = \"8+47=\";\ntable.getCell(6, 0).value = \"89-30=\";\ntable.getCell(6, 1).value = \"38+26=\";\ntable.getCell(6, 2).value = \"54+33=\";\ntable.getCell(6, 3).value = \"25-15=\";\ntable.getCell(6, 4).value = \"20+54=\";\ntable.getCell(7, 0).value = \"62-6=\";\ntable.getCell(7, 1).value = \"37-16=\";\ntable.getCell(7, 2).value = \"31-3=\";\ntable.getCell(7, 3).value = \"20+75=\";\ntable.getCell(7, 4).value = \"13+17=\";\ntable.getCell(8, 0).value = \"18+23=\";\ntable.getCell(8, 1).value = \"81-46=\";\ntable.getCell(8, 2).value = \"76-69=\";\ntable.getCell(8, 3).value = \"1+32=\";\ntable.getCell(8, 4).value = \"51+41=\";\ntable.getCell(9, 0).value = \"39-15=\";\ntable.getCell(9, 1).value = \"81-61=\";\ntable.getCell(9, 2).value = \"10+45=\";\ntable.getCell(9, 3).value = \"48+41=\";\ntable.getCell(9, 4).value = \"58+17=\";\ntable.getCell(10, 0).value = \"41+0=\";\ntable.getCell(10, 1).value = \"84-62=\";\ntable.getCell(10, 2).value = \"51-0=\";\ntable.getCell(10, 3).value = \"40-5=\";\ntable.getCell(10, 4).value = \"9+65=\";\ntable.getCell(11, 0).value = \"10+51=\";\ntable.getCell(11, 1).value = \"36+28=\";\ntable.getCell(11, 2).value = \"83+7=\";\ntable.getCell(11, 3).value = \"11+49=\";\ntable.getCell(11, 4).value = \"58+29=\";\ntable.getCell(12, 0).value = \"92-13=\";\ntable.getCell(12, 1).value = \"91-66=\";\ntable.getCell(12, 2).value = \"41-25=\";\ntable.getCell(12, 3).value = \"81+11=\";\ntable.getCell(12, 4).value = \"8+25=\";\ntable.getCell(13, 0).value = \"23+18=\";\ntable.getCell(13, 1).value = \"82-48=\";\ntable.getCell(13, 2).value = \"86-52=\";\ntable.getCell(13, 3).value = \"9+39=\";\ntable.getCell(13, 4).value = \"94-15=\";\ntable.getCell(14, 0).value = \"90+0=\";\ntable.getCell(14, 1).value = \"95-25=\";\ntable.getCell(14, 2).value = \"63-44=\";\ntable.getCell(14, 3).value = \"72+12=\";\ntable.getCell(14, 4).value = \"46+25=\";\ntable.getCell(15, 0).value = \"88-48=\";\ntable.getCell(15, 1).value = \"54-44=\";\ntable.getCell(15, 2).value = \"85-41=\";\ntable.getCell(15, 3).value = \"83+2=\";\ntable.getCell(15, 4).value = \"22-11=\";\ntable.getCell(16, 0).value = \"51-41=\";\ntable.getCell(16, 1).value = \"42+45=\";\ntable.getCell(16, 2).value = \"18+64=\";\ntable.getCell(16, 3).value = \"93-50=\";\ntable.getCell(16, 4).value = \"55+40=\";\ntable.getCell(17, 0).value = \"57-52=\";\ntable.getCell(17, 1).value = \"66-33=\";\ntable.getCell(17, 2).value = \"67-13=\";\ntable.getCell(17, 3).value = \"60+38=\";\ntable.getCell(17, 4).value = \"17+18=\";\ntable.getCell(18, 0).value = \"99-97=\";\ntable.getCell(18, 1).value = \"62-51=\";\ntable.getCell(18, 2).value = \"4+49=\";\ntable.getCell(18, 3).value = \"26+8=\";\ntable.getCell(18, 4).value = \"0+62=\";\ntable.getCell(19, 0).value = \"80-72=\";\ntable.getCell(19, 1).value = \"15-10=\";\ntable.getCell(19, 2).value = \"67+13=\";\ntable.getCell(19, 3).value = \"80-58=\";\ntable.getCell(19, 4).value = \"87+0=\";\n\nawait context.sync();\n", "ps1": "# The document's sole table is a 20x5 grid of arithmetic equations\n# (\"a+b=\" / \"a-b=\"). The commit replaces the text of every cell with a\n# newly generated equation while leaving all formatting untouched.\n# Address each cell positionally (1-based row/col, Word COM style) and\n# overwrite its Range.Text so run/paragraph formatting (font, size,\n# justification) is preserved exactly.\n$d = $word.ActiveDocument\n$tbl = $d.Tables.Item(1)\n\n$tbl.Cell(1,1).Range.Text = \"26+36=\"\n$tbl.Cell(1,2).Range.Text = \"64+4=\"\n$tbl.Cell(1,3).Range.Text = \"44+54=\"\n$tbl.Cell(1,4).Range.Text = \"98-12=\"\n$tbl.Cell(1,5).Range.Text = \"91-27=\"\n$tbl.Cell(2,1).Range.Text = \"44+24=\"\n$tbl.Cell(2,2).Range.Text = \"44-32=\"\n$tbl.Cell(2,3).Range.Text = \"37+18=\"\n$tbl.Cell(2,4).Range.Text = \"66+25=\"\n$tbl.Cell(2,5).Range.Text = \"4+28=\"\n$tbl.Cell(3,1).Range.Text = \"40+51=\"\n$tbl.Cell(3,2).Range.Text = \"13+28=\"\n$tbl.Cell(3,3).Range.Text = \"66+3=\"\n$tbl.Cell(3,4).Range.Text = \"21+29=\"\n$tbl.Cell(3,5).Range.Text = \"66-65=\"\n$tbl.Cell(4,1).Range.Text = \"42-42=\"\n$tbl.Cell(4,2).Range.Text = \"92-6=\"\n$tbl.Cell(4,3).Range.Text = \"1+77=\"\n$tbl.Cell(4,4).Range.Text = \"38+2=\"\n$tbl.Cell(4,5).Range.Text = \"8+75=\"\n$tbl.Cell(5,1).Range.Text = \"28-10=\"\n$tbl.Cell(5,2).Range.Text = \"98-9=\"\n$tbl.Cell(5,3).Range.Text = \"36+7=\"\n$tbl.Cell(5,4).Range.Text = \"33-22=\"\n$tbl.Cell(5,5).Range.Text = \"53+44=\"\n$tbl.Cell(6,1).Range.Text = \"2+61=\"\n$tbl.Cell(6,2).Range.Text = \"97-76=\"\n$tbl.Cell(6,3).Range.Text = \"5+93=\"\n$tbl.Cell(6,4).Range.Text = \"71-20=\"\n$tbl.Cell(6,5).Range.Text = \"8+47=\"\n$tbl.Cell(7,1).Range.Text = \"89-30=\"\n$tbl.Cell(7,2).Range.Text = \"38+26=\"\n$tbl.Cell(7,3).Range.Text = \"54+33=\"\n$tbl.Cell(7,4).Range.Text = \"25-15=\"\n$tbl.Cell(7,5).Range.Text = \"20+54=\"\n$tbl.Cell(8,1).Range.Text = \"62-6=\"\n$tbl.Cell(8,2).Range.Text = \"37-16=\"\n$tbl.Cell(8,3).Range.Text = \"31-3=\"\n$tbl.Cell(8,4).Range.Text = \"20+75=\"\n$tbl.Cell(8,5).Range.Text = \"13+17=\"\n$tbl.Cell(9,1).Range.Text = \"18+23=\"\n$tbl.Cell(9,2).Range.Text = \"81-46=\"\n$tbl.Cell(9,3).Range.Text = \"76-69=\"\n$tbl.Cell(9,4).Range.Text = \"1+32=\"\n$tbl.Cell(9,5).Range.Text = \"51+41=\"\n$tbl.Cell(10,1).Range.Text = \"39-15=\"\n$tbl.Cell(10,2).Range.Text = \"81-61=\"\n$tbl.Cell(10,3).Range.Text = \"10+45=\"\n$tbl.Cell(10,4).Range.Text = \"48+41=\"\n$tbl.Cell(10,5).Range.Text = \"58+17=\"\n$tbl.Cell(11,1).Range.Text = \"41+0=\"\n$tbl.Cell(11,2).Range.Text = \"84-62=\"\n$tbl.Cell(11,3).Range.Text = \"51-0=\"\n$tbl.Cell(11,4).Range.Text = \"40-5=\"\n$tbl.Cell(11,5).Range.Text = \"9+65=\"\n$tbl.Cell(12,1).Range.Text = \"10+51=\"\n$tbl.Cell(12,2).Range.Text = \"36+28=\"\n$tbl.Cell(12,3).Range.Text = \"83+7=\"\n$tbl.Cell(12,4).Range.Text = \"11+49=\"\n$tbl.Cell(12,5).Range.Text = \"58+29=\"\n$tbl.Cell(13,1).Range.Text = \"92-13=\"\n$tbl.Cell(13,2).Range.Text = \"91-66=\"\n$tbl.Cell(13,3).Range.Text = \"41-25=\"\n$tbl.Cell(13,4).Range.Text = \"81+11=\"\n$tbl.Cell(13,5).Range.Text = \"8+25=\"\n$tbl.Cell(14,1).Range.Text = \"23+18=\"\n$tbl.Cell(14,2).Range.Text = \"82-48=\"\n$tbl.Cell(14,3).Range.Text = \"86-52=\"\n$tbl.Cell(14,4).Range.Text = \"9+39=\"\n$tbl.Cell(14,5).Range.Text = \"94-15=\"\n$tbl.Cell(15,1).Range.Text = \"90+0=\"\n$tbl.Cell(15,2).Range.Text = \"95-25=\"\n$tbl.Cell(15,3).Range.Text = \"63-44=\"\n$tbl.Cell(15,4).Range.Text = \"72+12=\"\n$tbl.Cell(15,5).Range.Text = \"46+25=\"\n$tbl.Cell(16,1).Range.Text = \"88-48=\"\n$tbl.Cell(16,2).Range.Text = \"54-44=\"\n$tbl.Cell(16,3).Range.Text = \"85-41=\"\n$tbl.Cell(16,4).Range.Text = \"83+2=\"\n$tbl.Cell(16,5).Range.Text = \"22-11=\"\n$tbl.Cell(17,1).Range.Text = \"51-41=\"\n$tbl.Cell(17,2).Range.Text = \"42+45=\"\n$tbl.Cell(17,3).Range.Text = \"18+64=\"\n$tbl.Cell(17,4).Range.Text = \"93-50=\"\n$tbl.Cell(17,5).Range.Text = \"55+40=\"\n$tbl.Cell(18,1).Range.Text = \"57-52=\"\n$tbl.Cell(18,2).Range.Text = \"66-33=\"\n$tbl.Cell(18,3).Range.Text = \"67-13=\"\n$tbl.Cell(18,4).Range.Text = \"60+38=\"\n$tbl.Cell(18,5).Range.Text = \"17+18=\"\n$tbl.Cell(19,1).Range.Text = \"99-97=\"\n$tbl.Cell(19,2).Range.Text = \"62-51=\"\n$tbl.Cell(19,3).Range.Text = \"4+49=\"\n$tbl.Cell(19,4).Range.Text = \"26+8=\"\n$tbl.Cell(19,5).Range.Text = \"0+62=\"\n$tbl.Cell(20,1).Range.Text = \"80-72=\"\n$tbl.Cell(20,2).Range.Text = \"15-10=\"\n$tbl.Cell(20,3).Range.Text = \"67+13=\"\n$tbl.Cell(20,4).Range.Text = \"80-58=\"\n$tbl.Cell(20,5).Range.Text = \"87+0=\"\n"}
